# Natmi following Dr Hou advice
#
# The LR-pair table now covers 3 sending/target clusters (ECs, FAPs, sCs)
# instead of just FAPs/sCs, growing the data block from rows 2-5 to rows 2-10
# (dimension A1:T10). Row 1 (headers) is untouched. We rewrite every data
# cell (columns A-T) for rows 2-10 with the updated cluster labels and
# recomputed expression/specificity statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is the full A..T row content, in column order, for rows 2-10.
$rows = @(
  @("ECs",  "Fgf2", "Fgfr2", "ECs",  2, 0.6666666666666666, 0.6462393333333333, 1.938718,   0.03461850536298827, 0.03461850536298827, 1, 0.3333333333333333, 0.4871643333333333, 1.461493, 0.1914458468300136, 0.1914458468300136, 0.3148247539971111,  2.833422785974,     0.006627569075206657, 0.006627569075206655),
  @("ECs",  "Fgf2", "Fgfr2", "FAPs", 2, 0.6666666666666666, 0.6462393333333333, 1.938718,   0.03461850536298827, 0.03461850536298827, 3, 1,                  1.864751,            5.594253, 0.7328098752209857, 0.7328098752209857, 1.205075443072667,   10.845678987654,    0.02536878259538846,  0.02536878259538846),
  @("ECs",  "Fgf2", "Fgfr2", "sCs",  2, 0.6666666666666666, 0.6462393333333333, 1.938718,   0.03461850536298827, 0.03461850536298827, 3, 1,                  0.1927433333333333, 0.57823,  0.07574427794900063, 0.07574427794900063, 0.1245583232377778,  1.12102490914,      0.002622153692393153, 0.002622153692393152),
  @("FAPs", "Fgf2", "Fgfr2", "ECs",  3, 1,                  15.322826,           45.968478,  0.8208310864042159,  0.8208310864042158,  1, 0.3333333333333333, 0.4871643333333333, 1.461493, 0.1914458468300136, 0.1914458468300136, 7.464734313072666,   67.182608817654,    0.1571447024410552,   0.1571447024410551),
  @("FAPs", "Fgf2", "Fgfr2", "FAPs", 3, 1,                  15.322826,           45.968478,  0.8208310864042159,  0.8208310864042158,  3, 1,                  1.864751,            5.594253, 0.7328098752209857, 0.7328098752209857, 28.573255106326,     257.159295956934,   0.6015131260053795,   0.6015131260053795),
  @("FAPs", "Fgf2", "Fgfr2", "sCs",  3, 1,                  15.322826,           45.968478,  0.8208310864042159,  0.8208310864042158,  3, 1,                  0.1927433333333333, 0.57823,  0.07574427794900063, 0.07574427794900063, 2.953372559326667,   26.58035303394,     0.06217325795778108,  0.06217325795778107),
  @("sCs",  "Fgf2", "Fgfr2", "ECs",  3, 1,                  2.698388,            8.095164,   0.1445504082327959,  0.1445504082327959,  1, 0.3333333333333333, 0.4871643333333333, 1.461493, 0.1914458468300136, 0.1914458468300136, 1.314558391094667,   11.831025519852,    0.02767357531375178,  0.02767357531375178),
  @("sCs",  "Fgf2", "Fgfr2", "FAPs", 3, 1,                  2.698388,            8.095164,   0.1445504082327959,  0.1445504082327959,  3, 1,                  1.864751,            5.594253, 0.7328098752209857, 0.7328098752209857, 5.031821721388,      45.28639549249201,  0.1059279666202177,   0.1059279666202177),
  @("sCs",  "Fgf2", "Fgfr2", "sCs",  3, 1,                  2.698388,            8.095164,   0.1445504082327959,  0.1445504082327959,  3, 1,                  0.1927433333333333, 0.57823,  0.07574427794900063, 0.07574427794900063, 0.5200962977466667,  4.68086667972,      0.0109488662988264,   0.0109488662988264)
)

$firstRow = 2
for ($r = 0; $r -lt $rows.Length; $r++) {
  $rowVals = $rows[$r]
  $excelRow = $firstRow + $r
  for ($col = 0; $col -lt $rowVals.Length; $col++) {
    $ws.Cells.Item($excelRow, $col + 1).Value = $rowVals[$col]
  }
}

Write-Output "Updated $($rows.Length) data rows (rows 2-10, columns A-T)."
